$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2004.909  # H17 (was 2052.3125)
$ws.Cells.Item(17, 10).Value = 1721.8438  # J17 (was 1761.6451)
$ws.Cells.Item(17, 12).Value = 5165.5314  # L17 (was 5284.9353)
$ws.Cells.Item(17, 14).Value = -5501.5314  # N17 (was -5620.9353)
$ws.Cells.Item(32, 8).Value = 1500  # H32 (was 1250)
$ws.Cells.Item(32, 9).Value = 0  # I32 (was 1000)
$ws.Cells.Item(32, 11).Value = 0  # K32 (was 1000)
$ws.Cells.Item(32, 13).Value = $null  # M32 (was -674)
$ws.Cells.Item(33, 8).Value = 278.75  # H33 (was 270.1905)
$ws.Cells.Item(33, 10).Value = 117.666664  # J33 (was 115.8)
$ws.Cells.Item(33, 12).Value = 117.666664  # L33 (was 115.8)
$ws.Cells.Item(33, 14).Value = -575.666664  # N33 (was -573.8)
$ws.Cells.Item(86, 8).Value = 154886.38  # H86 (was 308905.5)
$ws.Cells.Item(86, 9).Value = 308850.25  # I86 (was 411544)
$ws.Cells.Item(86, 10).Value = 922.5  # J86 (was 990)
$ws.Cells.Item(86, 11).Value = 308850.25  # K86 (was 411544)
$ws.Cells.Item(86, 12).Value = 922.5  # L86 (was 990)
$ws.Cells.Item(86, 13).Value = -307727.25  # M86 (was -410421)
$ws.Cells.Item(86, 14).Value = -3168.5  # N86 (was -3236)
$ws.Cells.Item(89, 8).Value = 154886.38  # H89 (was 308905.5)
$ws.Cells.Item(89, 9).Value = 308850.25  # I89 (was 411544)
$ws.Cells.Item(89, 10).Value = 922.5  # J89 (was 990)
$ws.Cells.Item(89, 11).Value = 1544251.25  # K89 (was 2057720)
$ws.Cells.Item(89, 12).Value = 4612.5  # L89 (was 4950)
$ws.Cells.Item(89, 13).Value = -1538635.25  # M89 (was -2052104)
$ws.Cells.Item(89, 14).Value = -15844.5  # N89 (was -16182)
$ws.Cells.Item(92, 8).Value = 947474  # H92 (was 1231666.4)
$ws.Cells.Item(92, 9).Value = 1231516.5  # I92 (was 1759238.1)
$ws.Cells.Item(92, 11).Value = 1231516.5  # K92 (was 1759238.1)
$ws.Cells.Item(92, 13).Value = -1230268.5  # M92 (was -1757990.1)
$ws.Cells.Item(98, 8).Value = 2219.9443  # H98 (was 2144)
$ws.Cells.Item(98, 9).Value = 2363.9333  # I98 (was 2282)
$ws.Cells.Item(98, 11).Value = 2363.9333  # K98 (was 2282)
$ws.Cells.Item(98, 13).Value = -865.9333000000001  # M98 (was -784)
$ws.Cells.Item(122, 8).Value = 2219.9443  # H122 (was 2144)
$ws.Cells.Item(122, 9).Value = 2363.9333  # I122 (was 2282)
$ws.Cells.Item(122, 11).Value = 7091.7999  # K122 (was 6846)
$ws.Cells.Item(122, 13).Value = -4641.7999  # M122 (was -4396)
$ws.Cells.Item(132, 8).Value = 769.2461499999999  # H132 (was 729.59155)
$ws.Cells.Item(132, 9).Value = 713.2951  # I132 (was 677.7761)
$ws.Cells.Item(132, 10).Value = 1622.5  # J132 (was 1597.5)
$ws.Cells.Item(132, 11).Value = 2139.8853  # K132 (was 2033.3283)
$ws.Cells.Item(132, 12).Value = 4867.5  # L132 (was 4792.5)
$ws.Cells.Item(132, 13).Value = 390.1147000000001  # M132 (was 496.6716999999999)
$ws.Cells.Item(132, 14).Value = -9927.5  # N132 (was -9852.5)
$ws.Cells.Item(137, 8).Value = 2050.8572  # H137 (was 1522.2368)
$ws.Cells.Item(137, 9).Value = 1324  # I137 (was 1048.2174)
$ws.Cells.Item(137, 10).Value = 2596  # J137 (was 2249.0667)
$ws.Cells.Item(137, 11).Value = 3972  # K137 (was 3144.6522)
$ws.Cells.Item(137, 12).Value = 7788  # L137 (was 6747.2001)
$ws.Cells.Item(137, 13).Value = -1422  # M137 (was -594.6522)
$ws.Cells.Item(137, 14).Value = -12888  # N137 (was -11847.2001)
$ws.Cells.Item(138, 8).Value = 1774.2  # H138 (was 1827)
$ws.Cells.Item(138, 9).Value = 1238  # I138 (was 1218.5278)
$ws.Cells.Item(138, 10).Value = 2524.88  # J138 (was 2452.8572)
$ws.Cells.Item(138, 11).Value = 3714  # K138 (was 3655.5834)
$ws.Cells.Item(138, 12).Value = 7574.64  # L138 (was 7358.571599999999)
$ws.Cells.Item(138, 13).Value = 1426  # M138 (was 1484.4166)
$ws.Cells.Item(138, 14).Value = -17854.64  # N138 (was -17638.5716)
$ws.Cells.Item(140, 8).Value = 53445.26  # H140 (was 54516.668)
$ws.Cells.Item(140, 10).Value = 53445.26  # J140 (was 54516.668)
$ws.Cells.Item(140, 12).Value = 53445.26  # L140 (was 54516.668)
$ws.Cells.Item(140, 14).Value = -63805.26  # N140 (was -64876.668)

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2534.7087  # H32 (was 2648.3066)
$ws.Cells.Item(32, 9).Value = 2241.2026  # I32 (was 2346.1428)
$ws.Cells.Item(32, 11).Value = 2241.2026  # K32 (was 2346.1428)
$ws.Cells.Item(32, 13).Value = -1954.2026  # M32 (was -2059.1428)
$ws.Cells.Item(74, 8).Value = 1752.4445  # H74 (was 1815.3125)
$ws.Cells.Item(74, 9).Value = 1650.4  # I74 (was 1778.3334)
$ws.Cells.Item(74, 10).Value = 1880  # J74 (was 1862.8572)
$ws.Cells.Item(74, 11).Value = 1650.4  # K74 (was 1778.3334)
$ws.Cells.Item(74, 12).Value = 1880  # L74 (was 1862.8572)
$ws.Cells.Item(74, 13).Value = -776.4000000000001  # M74 (was -904.3334)
$ws.Cells.Item(74, 14).Value = -3628  # N74 (was -3610.8572)
$ws.Cells.Item(77, 8).Value = 1752.4445  # H77 (was 1815.3125)
$ws.Cells.Item(77, 9).Value = 1650.4  # I77 (was 1778.3334)
$ws.Cells.Item(77, 10).Value = 1880  # J77 (was 1862.8572)
$ws.Cells.Item(77, 11).Value = 8252  # K77 (was 8891.666999999999)
$ws.Cells.Item(77, 12).Value = 9400  # L77 (was 9314.286)
$ws.Cells.Item(77, 13).Value = -3884  # M77 (was -4523.666999999999)
$ws.Cells.Item(77, 14).Value = -18136  # N77 (was -18050.286)
$ws.Cells.Item(97, 8).Value = 1007.5  # H97 (was 1040.4546)
$ws.Cells.Item(97, 10).Value = 647.5  # J97 (was 650)
$ws.Cells.Item(97, 12).Value = 647.5  # L97 (was 650)
$ws.Cells.Item(97, 14).Value = -1639.5  # N97 (was -1642)
$ws.Cells.Item(102, 8).Value = 2787.1428  # H102 (was 0)
$ws.Cells.Item(102, 9).Value = 2787.1428  # I102 (was 0)
$ws.Cells.Item(102, 11).Value = 2787.1428  # K102 (was 0)
$ws.Cells.Item(102, 13).Value = -1165.1428  # M102 (was None)
$ws.Cells.Item(132, 8).Value = 1572.8937  # H132 (was 1624.4783)
$ws.Cells.Item(132, 9).Value = 952.0909  # I132 (was 969.34375)
$ws.Cells.Item(132, 10).Value = 3036.2144  # J132 (was 3121.9285)
$ws.Cells.Item(132, 11).Value = 2856.2727  # K132 (was 2908.03125)
$ws.Cells.Item(132, 12).Value = 9108.643199999999  # L132 (was 9365.7855)
$ws.Cells.Item(132, 13).Value = -326.2727  # M132 (was -378.03125)
$ws.Cells.Item(132, 14).Value = -14168.6432  # N132 (was -14425.7855)

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1856.1538  # H20 (was 1899.6)
$ws.Cells.Item(20, 9).Value = 1754.3043  # I20 (was 1799.0454)
$ws.Cells.Item(20, 11).Value = 1754.3043  # K20 (was 1799.0454)
$ws.Cells.Item(20, 13).Value = -1507.3043  # M20 (was -1552.0454)
$ws.Cells.Item(94, 8).Value = 840  # H94 (was 925)
$ws.Cells.Item(94, 9).Value = 866.6667  # I94 (was 1050)
$ws.Cells.Item(94, 11).Value = 866.6667  # K94 (was 1050)
$ws.Cells.Item(94, 13).Value = -415.6667  # M94 (was -599)
$ws.Cells.Item(105, 8).Value = 2461.077  # H105 (was 2416.2917)
$ws.Cells.Item(105, 9).Value = 2299.5  # I105 (was 2299.6)
$ws.Cells.Item(105, 10).Value = 2999.6667  # J105 (was 2999.75)
$ws.Cells.Item(105, 11).Value = 2299.5  # K105 (was 2299.6)
$ws.Cells.Item(105, 12).Value = 2999.6667  # L105 (was 2999.75)
$ws.Cells.Item(105, 13).Value = -552.5  # M105 (was -552.5999999999999)
$ws.Cells.Item(105, 14).Value = -6493.6667  # N105 (was -6493.75)
$ws.Cells.Item(134, 8).Value = 5012.4316  # H134 (was 4718.0215)
$ws.Cells.Item(134, 9).Value = 5739.3335  # I134 (was 5198.5674)
$ws.Cells.Item(134, 10).Value = 2831.7273  # J134 (was 2940)
$ws.Cells.Item(134, 11).Value = 17218.0005  # K134 (was 15595.7022)
$ws.Cells.Item(134, 12).Value = 8495.1819  # L134 (was 8820)
$ws.Cells.Item(134, 13).Value = -14683.0005  # M134 (was -13060.7022)
$ws.Cells.Item(134, 14).Value = -13565.1819  # N134 (was -13890)
$ws.Cells.Item(138, 8).Value = 0  # H138 (was 40000)
$ws.Cells.Item(138, 10).Value = 0  # J138 (was 40000)
$ws.Cells.Item(138, 12).Value = 0  # L138 (was 40000)
$ws.Cells.Item(138, 14).Value = $null  # N138 (was -50280)

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 129.6923  # H7 (was 215.5)
$ws.Cells.Item(7, 9).Value = 48.714287  # I7 (was 77.25)
$ws.Cells.Item(7, 10).Value = 224.16667  # J7 (was 353.75)
$ws.Cells.Item(7, 11).Value = 48.714287  # K7 (was 77.25)
$ws.Cells.Item(7, 12).Value = 224.16667  # L7 (was 353.75)
$ws.Cells.Item(7, 13).Value = 64.285713  # M7 (was 35.75)
$ws.Cells.Item(7, 14).Value = -450.16667  # N7 (was -579.75)
$ws.Cells.Item(31, 8).Value = 2603.3333  # H31 (was 2186.0952)
$ws.Cells.Item(31, 9).Value = 2213.6  # I31 (was 2080.4546)
$ws.Cells.Item(31, 10).Value = 3382.8  # J31 (was 2302.3)
$ws.Cells.Item(31, 11).Value = 2213.6  # K31 (was 2080.4546)
$ws.Cells.Item(31, 12).Value = 3382.8  # L31 (was 2302.3)
$ws.Cells.Item(31, 13).Value = -1918.6  # M31 (was -1785.4546)
$ws.Cells.Item(31, 14).Value = -3972.8  # N31 (was -2892.3)
$ws.Cells.Item(34, 8).Value = 2603.3333  # H34 (was 2186.0952)
$ws.Cells.Item(34, 9).Value = 2213.6  # I34 (was 2080.4546)
$ws.Cells.Item(34, 10).Value = 3382.8  # J34 (was 2302.3)
$ws.Cells.Item(34, 11).Value = 2213.6  # K34 (was 2080.4546)
$ws.Cells.Item(34, 12).Value = 3382.8  # L34 (was 2302.3)
$ws.Cells.Item(34, 13).Value = -2011.6  # M34 (was -1878.4546)
$ws.Cells.Item(34, 14).Value = -3786.8  # N34 (was -2706.3)
$ws.Cells.Item(58, 8).Value = 1501335.8  # H58 (was 1451334.5)
$ws.Cells.Item(58, 9).Value = 1891900.4  # I58 (was 1813125.4)
$ws.Cells.Item(58, 11).Value = 1891900.4  # K58 (was 1813125.4)
$ws.Cells.Item(58, 13).Value = -1891697.4  # M58 (was -1812922.4)
$ws.Cells.Item(94, 8).Value = 1018.4  # H94 (was 1014.9)
$ws.Cells.Item(94, 9).Value = 999.3333  # I94 (was 809.25)
$ws.Cells.Item(94, 10).Value = 1026.5714  # J94 (was 1152)
$ws.Cells.Item(94, 11).Value = 999.3333  # K94 (was 809.25)
$ws.Cells.Item(94, 12).Value = 1026.5714  # L94 (was 1152)
$ws.Cells.Item(94, 13).Value = -548.3333  # M94 (was -358.25)
$ws.Cells.Item(94, 14).Value = -1928.5714  # N94 (was -2054)
$ws.Cells.Item(99, 8).Value = 2099.7778  # H99 (was 2137.25)
$ws.Cells.Item(99, 9).Value = 2056.8572  # I99 (was 2085.4285)
$ws.Cells.Item(99, 10).Value = 2250  # J99 (was 2500)
$ws.Cells.Item(99, 11).Value = 2056.8572  # K99 (was 2085.4285)
$ws.Cells.Item(99, 12).Value = 2250  # L99 (was 2500)
$ws.Cells.Item(99, 13).Value = -558.8571999999999  # M99 (was -587.4285)
$ws.Cells.Item(99, 14).Value = -5246  # N99 (was -5496)
$ws.Cells.Item(126, 8).Value = 2099.7778  # H126 (was 2137.25)
$ws.Cells.Item(126, 9).Value = 2056.8572  # I126 (was 2085.4285)
$ws.Cells.Item(126, 10).Value = 2250  # J126 (was 2500)
$ws.Cells.Item(126, 11).Value = 6170.571599999999  # K126 (was 6256.2855)
$ws.Cells.Item(126, 12).Value = 6750  # L126 (was 7500)
$ws.Cells.Item(126, 13).Value = -3700.571599999999  # M126 (was -3786.2855)
$ws.Cells.Item(126, 14).Value = -11690  # N126 (was -12440)
$ws.Cells.Item(132, 8).Value = 2133.3057  # H132 (was 2169.9714)
$ws.Cells.Item(132, 9).Value = 1352.75  # I132 (was 1379.2106)
$ws.Cells.Item(132, 11).Value = 4058.25  # K132 (was 4137.6318)
$ws.Cells.Item(132, 13).Value = -1528.25  # M132 (was -1607.6318)
$ws.Cells.Item(134, 8).Value = 1591.55  # H134 (was 1603.2307)
$ws.Cells.Item(134, 9).Value = 1393.4722  # I134 (was 1431.3611)
$ws.Cells.Item(134, 10).Value = 3374.25  # J134 (was 3665.6667)
$ws.Cells.Item(134, 11).Value = 4180.4166  # K134 (was 4294.0833)
$ws.Cells.Item(134, 12).Value = 10122.75  # L134 (was 10997.0001)
$ws.Cells.Item(134, 13).Value = -1645.4166  # M134 (was -1759.0833)
$ws.Cells.Item(134, 14).Value = -15192.75  # N134 (was -16067.0001)
$ws.Cells.Item(136, 8).Value = 1501335.8  # H136 (was 1451334.5)
$ws.Cells.Item(136, 9).Value = 1891900.4  # I136 (was 1813125.4)
$ws.Cells.Item(136, 11).Value = 5675701.199999999  # K136 (was 5439376.199999999)
$ws.Cells.Item(136, 13).Value = -5673151.199999999  # M136 (was -5436826.199999999)

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 1999  # H75 (was 1006.5)
$ws.Cells.Item(75, 9).Value = 0  # I75 (was 13)
$ws.Cells.Item(75, 10).Value = 1999  # J75 (was 2000)
$ws.Cells.Item(75, 11).Value = 0  # K75 (was 39)
$ws.Cells.Item(75, 12).Value = 5997  # L75 (was 6000)
$ws.Cells.Item(75, 13).Value = $null  # M75 (was 959)
$ws.Cells.Item(75, 14).Value = -7993  # N75 (was -7996)
$ws.Cells.Item(78, 8).Value = 1999  # H78 (was 1006.5)
$ws.Cells.Item(78, 9).Value = 0  # I78 (was 13)
$ws.Cells.Item(78, 10).Value = 1999  # J78 (was 2000)
$ws.Cells.Item(78, 11).Value = 0  # K78 (was 117)
$ws.Cells.Item(78, 12).Value = 17991  # L78 (was 18000)
$ws.Cells.Item(78, 13).Value = $null  # M78 (was 4875)
$ws.Cells.Item(78, 14).Value = -27975  # N78 (was -27984)
$ws.Cells.Item(107, 8).Value = 890.5625  # H107 (was 968.6429000000001)
$ws.Cells.Item(107, 9).Value = 637.6  # I107 (was 833.3333)
$ws.Cells.Item(107, 11).Value = 1912.8  # K107 (was 2499.9999)
$ws.Cells.Item(107, 13).Value = 7.199999999999818  # M107 (was -579.9998999999998)
$ws.Cells.Item(113, 8).Value = 9118.083000000001  # H113 (was 6969)
$ws.Cells.Item(113, 9).Value = 50500.5  # I113 (was 25500.25)
$ws.Cells.Item(113, 10).Value = 841.6  # J113 (was 791.9167)
$ws.Cells.Item(113, 11).Value = 151501.5  # K113 (was 76500.75)
$ws.Cells.Item(113, 12).Value = 2524.8  # L113 (was 2375.7501)
$ws.Cells.Item(113, 13).Value = -149331.5  # M113 (was -74330.75)
$ws.Cells.Item(113, 14).Value = -6864.8  # N113 (was -6715.7501)
$ws.Cells.Item(132, 8).Value = 1152.4166  # H132 (was 1099.5714)
$ws.Cells.Item(132, 9).Value = 0  # I132 (was 780)
$ws.Cells.Item(132, 10).Value = 1152.4166  # J132 (was 1152.8334)
$ws.Cells.Item(132, 11).Value = 0  # K132 (was 7020)
$ws.Cells.Item(132, 12).Value = 10371.7494  # L132 (was 10375.5006)
$ws.Cells.Item(132, 13).Value = $null  # M132 (was -4490)
$ws.Cells.Item(132, 14).Value = -15431.7494  # N132 (was -15435.5006)
$ws.Cells.Item(137, 8).Value = 4720.3  # H137 (was 5424.125)
$ws.Cells.Item(137, 9).Value = 3734  # I137 (was 4090)
$ws.Cells.Item(137, 10).Value = 5706.6  # J137 (was 6758.25)
$ws.Cells.Item(137, 11).Value = 11202  # K137 (was 12270)
$ws.Cells.Item(137, 12).Value = 17119.8  # L137 (was 20274.75)
$ws.Cells.Item(137, 13).Value = -6102  # M137 (was -7170)
$ws.Cells.Item(137, 14).Value = -27319.8  # N137 (was -30474.75)
$ws.Cells.Item(140, 8).Value = 1922.6666  # H140 (was 1912.9445)
$ws.Cells.Item(140, 9).Value = 1401.7273  # I140 (was 1369.909)
$ws.Cells.Item(140, 11).Value = 4205.1819  # K140 (was 4109.727000000001)
$ws.Cells.Item(140, 13).Value = 974.8181000000004  # M140 (was 1070.272999999999)

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5333.3335  # H70 (was 4880)
$ws.Cells.Item(70, 9).Value = 7000  # I70 (was 5450)
$ws.Cells.Item(70, 11).Value = 7000  # K70 (was 5450)
$ws.Cells.Item(70, 13).Value = -6730  # M70 (was -5180)
$ws.Cells.Item(73, 8).Value = 5333.3335  # H73 (was 4880)
$ws.Cells.Item(73, 9).Value = 7000  # I73 (was 5450)
$ws.Cells.Item(73, 11).Value = 7000  # K73 (was 5450)
$ws.Cells.Item(73, 13).Value = -6064  # M73 (was -4514)
$ws.Cells.Item(80, 8).Value = 1653.7778  # H80 (was 1798.125)
$ws.Cells.Item(80, 9).Value = 1798  # I80 (was 1955)
$ws.Cells.Item(80, 10).Value = 500  # J80 (was 700)
$ws.Cells.Item(80, 11).Value = 1798  # K80 (was 1955)
$ws.Cells.Item(80, 12).Value = 500  # L80 (was 700)
$ws.Cells.Item(80, 13).Value = -800  # M80 (was -957)
$ws.Cells.Item(80, 14).Value = -2496  # N80 (was -2696)
$ws.Cells.Item(83, 8).Value = 1653.7778  # H83 (was 1798.125)
$ws.Cells.Item(83, 9).Value = 1798  # I83 (was 1955)
$ws.Cells.Item(83, 10).Value = 500  # J83 (was 700)
$ws.Cells.Item(83, 11).Value = 8990  # K83 (was 9775)
$ws.Cells.Item(83, 12).Value = 2500  # L83 (was 3500)
$ws.Cells.Item(83, 13).Value = -3998  # M83 (was -4783)
$ws.Cells.Item(83, 14).Value = -12484  # N83 (was -13484)
$ws.Cells.Item(102, 8).Value = 1844.6888  # H102 (was 1883.841)
$ws.Cells.Item(102, 9).Value = 1815.1111  # I102 (was 1863.4857)
$ws.Cells.Item(102, 11).Value = 1815.1111  # K102 (was 1863.4857)
$ws.Cells.Item(102, 13).Value = -193.1111000000001  # M102 (was -241.4857)
$ws.Cells.Item(113, 8).Value = 1224.75  # H113 (was 1325)
$ws.Cells.Item(113, 9).Value = 999  # I113 (was 0)
$ws.Cells.Item(113, 10).Value = 1300  # J113 (was 1325)
$ws.Cells.Item(113, 11).Value = 999  # K113 (was 0)
$ws.Cells.Item(113, 12).Value = 1300  # L113 (was 1325)
$ws.Cells.Item(113, 13).Value = 1171  # M113 (was None)
$ws.Cells.Item(113, 14).Value = -5640  # N113 (was -5665)
$ws.Cells.Item(122, 8).Value = 2500  # H122 (was 2366.6667)
$ws.Cells.Item(122, 9).Value = 0  # I122 (was 2100)
$ws.Cells.Item(122, 11).Value = 0  # K122 (was 6300)
$ws.Cells.Item(122, 13).Value = $null  # M122 (was -3850)
$ws.Cells.Item(136, 8).Value = 10990.1  # H136 (was 10602.454)
$ws.Cells.Item(136, 10).Value = 10990.1  # J136 (was 10602.454)
$ws.Cells.Item(136, 12).Value = 32970.3  # L136 (was 31807.362)
$ws.Cells.Item(136, 14).Value = -38070.3  # N136 (was -36907.362)

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 13928.728  # H40 (was 15054.7)
$ws.Cells.Item(40, 9).Value = 27093.25  # I40 (was 35234.668)
$ws.Cells.Item(40, 11).Value = 27093.25  # K40 (was 35234.668)
$ws.Cells.Item(40, 13).Value = -26957.25  # M40 (was -35098.668)
$ws.Cells.Item(43, 8).Value = 10455.333  # H43 (was 10808.4)
$ws.Cells.Item(43, 10).Value = 10455.333  # J43 (was 10808.4)
$ws.Cells.Item(43, 12).Value = 10455.333  # L43 (was 10808.4)
$ws.Cells.Item(43, 14).Value = -10841.333  # N43 (was -11194.4)
$ws.Cells.Item(61, 8).Value = 2484.1538  # H61 (was 1657.625)
$ws.Cells.Item(61, 9).Value = 2221.5557  # I61 (was 1607.9375)
$ws.Cells.Item(61, 10).Value = 3075  # J61 (was 1757)
$ws.Cells.Item(61, 11).Value = 2221.5557  # K61 (was 1607.9375)
$ws.Cells.Item(61, 12).Value = 3075  # L61 (was 1757)
$ws.Cells.Item(61, 13).Value = -2019.5557  # M61 (was -1405.9375)
$ws.Cells.Item(61, 14).Value = -3479  # N61 (was -2161)
$ws.Cells.Item(68, 8).Value = 2119.5557  # H68 (was 2264.625)
$ws.Cells.Item(68, 9).Value = 1679.5  # I68 (was 1852.8334)
$ws.Cells.Item(68, 10).Value = 2999.6667  # J68 (was 3500)
$ws.Cells.Item(68, 11).Value = 1679.5  # K68 (was 1852.8334)
$ws.Cells.Item(68, 12).Value = 2999.6667  # L68 (was 3500)
$ws.Cells.Item(68, 13).Value = -930.5  # M68 (was -1103.8334)
$ws.Cells.Item(68, 14).Value = -4497.6667  # N68 (was -4998)
$ws.Cells.Item(71, 8).Value = 2119.5557  # H71 (was 2264.625)
$ws.Cells.Item(71, 9).Value = 1679.5  # I71 (was 1852.8334)
$ws.Cells.Item(71, 10).Value = 2999.6667  # J71 (was 3500)
$ws.Cells.Item(71, 11).Value = 8397.5  # K71 (was 9264.166999999999)
$ws.Cells.Item(71, 12).Value = 14998.3335  # L71 (was 17500)
$ws.Cells.Item(71, 13).Value = -4653.5  # M71 (was -5520.166999999999)
$ws.Cells.Item(71, 14).Value = -22486.3335  # N71 (was -24988)
$ws.Cells.Item(113, 8).Value = 2484.1538  # H113 (was 1657.625)
$ws.Cells.Item(113, 9).Value = 2221.5557  # I113 (was 1607.9375)
$ws.Cells.Item(113, 10).Value = 3075  # J113 (was 1757)
$ws.Cells.Item(113, 11).Value = 2221.5557  # K113 (was 1607.9375)
$ws.Cells.Item(113, 12).Value = 3075  # L113 (was 1757)
$ws.Cells.Item(113, 13).Value = -51.55569999999989  # M113 (was 562.0625)
$ws.Cells.Item(113, 14).Value = -7415  # N113 (was -6097)
$ws.Cells.Item(136, 8).Value = 3145.037  # H136 (was 3503.1738)
$ws.Cells.Item(136, 9).Value = 1740.8334  # I136 (was 1928)
$ws.Cells.Item(136, 11).Value = 5222.5002  # K136 (was 5784)
$ws.Cells.Item(136, 13).Value = -2672.5002  # M136 (was -3234)

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 0  # H4 (was 6000)
$ws.Cells.Item(4, 10).Value = 0  # J4 (was 6000)
$ws.Cells.Item(4, 12).Value = 0  # L4 (was 6000)
$ws.Cells.Item(4, 14).Value = $null  # N4 (was -6226)
$ws.Cells.Item(70, 8).Value = 24500  # H70 (was 29777.5)
$ws.Cells.Item(70, 9).Value = 8000  # I70 (was 0)
$ws.Cells.Item(70, 10).Value = 30000  # J70 (was 29777.5)
$ws.Cells.Item(70, 11).Value = 8000  # K70 (was 0)
$ws.Cells.Item(70, 12).Value = 30000  # L70 (was 29777.5)
$ws.Cells.Item(70, 13).Value = -7685  # M70 (was None)
$ws.Cells.Item(70, 14).Value = -30630  # N70 (was -30407.5)
$ws.Cells.Item(73, 8).Value = 24500  # H73 (was 29777.5)
$ws.Cells.Item(73, 9).Value = 8000  # I73 (was 0)
$ws.Cells.Item(73, 10).Value = 30000  # J73 (was 29777.5)
$ws.Cells.Item(73, 11).Value = 8000  # K73 (was 0)
$ws.Cells.Item(73, 12).Value = 30000  # L73 (was 29777.5)
$ws.Cells.Item(73, 13).Value = -6908  # M73 (was None)
$ws.Cells.Item(73, 14).Value = -32184  # N73 (was -31961.5)
$ws.Cells.Item(126, 8).Value = 4199.8125  # H126 (was 4213.1875)
$ws.Cells.Item(126, 9).Value = 2067.875  # I126 (was 2751.4)
$ws.Cells.Item(126, 10).Value = 6331.75  # J126 (was 4877.636)
$ws.Cells.Item(126, 11).Value = 6203.625  # K126 (was 8254.200000000001)
$ws.Cells.Item(126, 12).Value = 18995.25  # L126 (was 14632.908)
$ws.Cells.Item(126, 13).Value = -3733.625  # M126 (was -5784.200000000001)
$ws.Cells.Item(126, 14).Value = -23935.25  # N126 (was -19572.908)
$ws.Cells.Item(136, 8).Value = 12347218  # H136 (was 12347223)
$ws.Cells.Item(136, 9).Value = 15874303  # I136 (was 16341189)
$ws.Cells.Item(136, 10).Value = 2418.7  # J136 (was 2234.182)
$ws.Cells.Item(136, 11).Value = 47622909  # K136 (was 49023567)
$ws.Cells.Item(136, 12).Value = 7256.099999999999  # L136 (was 6702.545999999999)
$ws.Cells.Item(136, 13).Value = -47620359  # M136 (was -49021017)
$ws.Cells.Item(136, 14).Value = -12356.1  # N136 (was -11802.546)
